# Trade #142 (global trade #170) closed at 2026-02-18 00:41:59 - unknown UNKNOWN +0.000%
# Also records two brand-new OPEN trades (#199 MarketMaking, #200 EMAArbitrage)
# and refreshes the roll-up Summary / Strategy Status numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal date-like / time-like string into a cell without
# letting Excel's automatic type inference turn it into a date/time serial.
# ---------------------------------------------------------------------------
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) Summary sheet roll-up numbers
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 170     # Total Trades
$summary.Range("B9").Value = 45.29   # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 66     # Trades
$status.Range("G6").Value = 46.97  # Win Rate %

# ---------------------------------------------------------------------------
# 3) All Trades sheet - close trade #170 (row 171) and append two new rows
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out existing open trade (row 171, Trade # 170, MarketMaking)
$allTrades.Range("G171").Value = 0.88
$allTrades.Range("H171").Value = "CLOSED"
$allTrades.Range("K171").Value = 99.34
$allTrades.Range("L171").Value = "early_exit"
$allTrades.Range("M171").Value = 0.19

# New row 200: Trade # 199 (MarketMaking, OPEN)
$allTrades.Range("A200").Value = 199
Set-TextCell $allTrades "B200" "2026-02-18"
Set-TextCell $allTrades "C200" "00:41:52"
$allTrades.Range("D200").Value = "MarketMaking"
$allTrades.Range("E200").Value = "DOWN"
$allTrades.Range("F200").Value = 0.88
$allTrades.Range("H200").Value = "OPEN"
$allTrades.Range("I200").Value = 0
$allTrades.Range("J200").Value = 0
$allTrades.Range("K200").Value = 99.33858346467945
$allTrades.Range("M200").Value = 0
$allTrades.Range("N200").Value = 0
$allTrades.Range("O200").Value = 0
$allTrades.Range("P200").Value = 0.6
$allTrades.Range("Q200").Value = "Normal spread capture: 198 bps"

# New row 201: Trade # 200 (EMAArbitrage, OPEN)
$allTrades.Range("A201").Value = 200
Set-TextCell $allTrades "B201" "2026-02-18"
Set-TextCell $allTrades "C201" "00:41:53"
$allTrades.Range("D201").Value = "EMAArbitrage"
$allTrades.Range("E201").Value = "DOWN"
$allTrades.Range("F201").Value = 0.89
$allTrades.Range("H201").Value = "OPEN"
$allTrades.Range("I201").Value = 0
$allTrades.Range("J201").Value = 0
$allTrades.Range("K201").Value = 100.270616878256
$allTrades.Range("M201").Value = 0
$allTrades.Range("N201").Value = 0
$allTrades.Range("O201").Value = 0
$allTrades.Range("P201").Value = 0.9
$allTrades.Range("Q201").Value = "EMA:down, RSI:50.0, ROC:-45.11% | 2/3 DOWN"

# ---------------------------------------------------------------------------
# 4) MarketMaking sheet - close trade #170 (row 67) and append new row 85
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close out existing open trade (row 67, Trade # 170)
$mm.Range("G67").Value = 0.88
$mm.Range("H67").Value = "CLOSED"
$mm.Range("K67").Value = 99.34
$mm.Range("P67").Value = "early_exit"
$mm.Range("Q67").Value = 0.19

# New row 85: Trade # 199 (OPEN)
$mm.Range("A85").Value = 199
Set-TextCell $mm "B85" "2026-02-18"
Set-TextCell $mm "C85" "00:41:52"
$mm.Range("D85").Value = "MarketMaking"
$mm.Range("E85").Value = "DOWN"
$mm.Range("F85").Value = 0.88
$mm.Range("H85").Value = "OPEN"
$mm.Range("I85").Value = 0
$mm.Range("J85").Value = 0
$mm.Range("K85").Value = 99.33858346467945
$mm.Range("L85").Value = 0
$mm.Range("M85").Value = 0
$mm.Range("N85").Value = 0.6
$mm.Range("O85").Value = "Normal spread capture: 198 bps"
$mm.Range("Q85").Value = 0

# ---------------------------------------------------------------------------
# 5) EMAArbitrage sheet - append new row 11
# ---------------------------------------------------------------------------
$ema = $wb.Worksheets.Item("EMAArbitrage")

# New row 11: Trade # 200 (OPEN)
$ema.Range("A11").Value = 200
Set-TextCell $ema "B11" "2026-02-18"
Set-TextCell $ema "C11" "00:41:53"
$ema.Range("D11").Value = "EMAArbitrage"
$ema.Range("E11").Value = "DOWN"
$ema.Range("F11").Value = 0.89
$ema.Range("H11").Value = "OPEN"
$ema.Range("I11").Value = 0
$ema.Range("J11").Value = 0
$ema.Range("K11").Value = 100.270616878256
$ema.Range("L11").Value = 0
$ema.Range("M11").Value = 0
$ema.Range("N11").Value = 0.9
$ema.Range("O11").Value = "EMA:down, RSI:50.0, ROC:-45.11% | 2/3 DOWN"
$ema.Range("Q11").Value = 0

Write-Output "edit applied"
